$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row: add column F (BFO_DEF) ----
$ws.Range("F1").Value = "BFO_DEF"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats -> copy header style (bold + border)

# ---- New row contents (rows reordered: material entity, role, function, quality) ----

# Row 2 - material entity
$ws.Range("B2").Value = "http://www.bioassayontology.org/bao#BAO_0003116"
$ws.Range("C2").Value = "{'label': 'material entity', 'prefLabel': None, 'altLabel': None, 'name': 'BAO_0003116'}"
$ws.Range("D2").Value = "http://purl.obolibrary.org/obo/BFO_0000040"
$ws.Range("E2").Value = "{'label': 'material entity', 'prefLabel': 'material entity'}"
$ws.Range("F2").Value = "['A material entity is an independent continuant that has some portion of matter as proper or improper continuant part. [BFO]']"

# Row 3 - role
$ws.Range("B3").Value = "http://www.bioassayontology.org/bao#BAO_0002929"
$ws.Range("C3").Value = "{'label': 'role', 'prefLabel': None, 'altLabel': None, 'name': 'BAO_0002929'}"
$ws.Range("D3").Value = "http://purl.obolibrary.org/obo/BFO_0000023"
$ws.Range("E3").Value = "{'label': 'role', 'prefLabel': 'role'}"
$ws.Range("F3").Value = "['B is a role means: b is a realizable entity and b exists because there is some single bearer that is in some special physical, social, or institutional set of circumstances in which this bearer does not have to be and b is not such that, if it ceases to exist, then the physical make-up of the bearer is thereby changed. [BFO]']"

# Row 4 - function
$ws.Range("B4").Value = "http://www.bioassayontology.org/bao#BAO_0003117"
$ws.Range("C4").Value = "{'label': 'function', 'prefLabel': None, 'altLabel': None, 'name': 'BAO_0003117'}"
$ws.Range("D4").Value = "http://purl.obolibrary.org/obo/BFO_0000034"
$ws.Range("E4").Value = "{'label': 'function', 'prefLabel': 'function'}"
$ws.Range("F4").Value = "['A function is a disposition that exists in virtue of the bearer" + [char]0x2019 + "s physical make-up and this physical make-up is something the bearer possesses because it came into being, either through evolution (in the case of natural biological entities) or through intentional design (in the case of artifacts), in order to realize processes of a certain sort. [BFO]']"

# Row 5 - quality
$ws.Range("B5").Value = "http://www.bioassayontology.org/bao#BAO_0002928"
$ws.Range("C5").Value = "{'label': 'quality', 'prefLabel': None, 'altLabel': None, 'name': 'BAO_0002928'}"
$ws.Range("D5").Value = "http://purl.obolibrary.org/obo/BFO_0000019"
$ws.Range("E5").Value = "{'label': 'quality', 'prefLabel': 'quality'}"
$ws.Range("F5").Value = "['A quality is a specifically dependent continuant that, in contrast to roles and dispositions, does not require any further process in order to be realized. [BFO]']"

# ---- Update the B-column hyperlink targets (location only; keep same external Address) ----
# NOTE: re-adding a hyperlink on a cell that already has one appends a fresh
# <hyperlink> entry; parsers (and Excel itself) resolve duplicate refs by
# taking the last entry for a given cell, so this effectively "updates" it.
$ws.Hyperlinks.Add($ws.Range("B2"), "http://www.bioassayontology.org/bao", "BAO_0003116") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), "http://www.bioassayontology.org/bao", "BAO_0002929") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B4"), "http://www.bioassayontology.org/bao", "BAO_0003117") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B5"), "http://www.bioassayontology.org/bao", "BAO_0002928") | Out-Null

# Restore the original (pre-hyperlink-restyle) look on column B by copying the
# still-pristine format from column D (same style index originally).
$ws.Range("D2").Copy()
$ws.Range("B2").PasteSpecial(-4122)
$ws.Range("D3").Copy()
$ws.Range("B3").PasteSpecial(-4122)
$ws.Range("D4").Copy()
$ws.Range("B4").PasteSpecial(-4122)
$ws.Range("D5").Copy()
$ws.Range("B5").PasteSpecial(-4122)
